$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns keep their literal text representation (avoid numeric coercion)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.109.19'
$ws.Range('E2').Value = '  +0.46%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.916.52'
$ws.Range('E3').Value = '  +2.44%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.14'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5069'
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4069'
$ws.Range('E8').Value = '  +3.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08333'
$ws.Range('E9').Value = '  +1.53%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.116'
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.03'
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '24.19'
$ws.Range('E12').Value = '  +5.58%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.419'
$ws.Range('E13').Value = '  +2.10%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.915.57'
$ws.Range('E14').Value = '  +2.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.250'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('E17').Value = '  +0.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001096'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06514'
$ws.Range('E19').Value = '  +2.24%  '
$ws.Range('E20').Value = '  +3.40%  '
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.950'
$ws.Range('E22').Value = '  +2.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '30.114.75'
$ws.Range('E23').Value = '  +0.47%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.35'
$ws.Range('E24').Value = '  +2.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.195'
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.137.19'
$ws.Range('E26').Value = '  +2.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '21.83'
$ws.Range('E27').Value = '  +4.20%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '162.70'
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.255'
$ws.Range('E29').Value = '  +0.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.81'
$ws.Range('E30').Value = '  +1.10%  '
$ws.Range('E31').Value = '  +6.30%  '
$ws.Range('E32').Value = '  +1.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.938'
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.786'
$ws.Range('E34').Value = '  +1.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02445'
$ws.Range('E35').Value = '  +0.39%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.304'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06431'
$ws.Range('E37').Value = '  +0.95%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.212'
$ws.Range('E38').Value = '  +3.14%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.2143'
$ws.Range('E39').Value = '  +0.11%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6479'
$ws.Range('E40').Value = '  +2.78%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.577'
$ws.Range('E41').Value = '  +0.76%  '
$ws.Range('E42').Value = '  +1.51%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.212'
$ws.Range('E43').Value = '  +0.38%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.35'
$ws.Range('E44').Value = '  +3.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.183'
$ws.Range('E45').Value = '  +8.85%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6042'
$ws.Range('E46').Value = '  +2.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.623'
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '122.40'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.209'
$ws.Range('E49').Value = '  +0.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.137'
$ws.Range('E50').Value = '  +1.49%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06832'
$ws.Range('E51').Value = '  +1.05%  '
